# Insert two new achievement bullet paragraphs into the "KEY ACHIEVEMENTS
# AND IMPACT" section, right after the existing "Expert methodology
# validated at highest judicial level" bullet and before the
# "TECHNICAL SKILLS" heading.

$d = $word.ActiveDocument

# Locate the anchor paragraph via Find so we don't depend on a hard-coded
# paragraph index, then resolve it to a real Paragraph object/index.
$found = $d.Content.Duplicate
$ok = $found.Find.Execute("Expert methodology validated at highest judicial level", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorIndex = $found.Paragraphs.Item(1).Index

# Work off of the real paragraph's Range (includes its end-of-paragraph
# mark), collapsed to its end, so InsertParagraphAfter lands between this
# paragraph and the next one ("TECHNICAL SKILLS").
$anchorPara = $d.Paragraphs.Item($anchorIndex)
$rng = $anchorPara.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$rng.InsertParagraphAfter()

# First new paragraph: plain bullet text.
$p1 = $d.Paragraphs.Item($anchorIndex + 1)
$p1.Range.Text = "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"

# Second new paragraph: "• " + bold/colored "178%" + trailing text.
$p2 = $d.Paragraphs.Item($anchorIndex + 2)
$p2.Range.Text = "• 178% accuracy improvement in racial classification algorithms"

# Bold + color the "178%" run to match the styling used elsewhere in the
# document for similar inline stats (e.g. "73.5%", "$4.7M").
$statRange = $p2.Range.Duplicate
$statRange.Find.Execute("178%", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$statRange.Bold = 1
$statRange.Font.Color = 5258796
